# Insert a new "Exceptional items" column into the Quarterly sheet,
# between "P/l before exceptional items & tax" (col K) and the old
# "P/l before tax" (old col L). Everything from the old column L onward
# shifts one column to the right (L->M, M->N, ... T->U).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Insert a whole new column before column L; this shifts existing data
# (including formatting) in columns L:T to M:U, and leaves the new
# column L blank.
$ws.Columns("L:L").Insert()

# Populate the header cells for the newly inserted column.
$ws.Range("L1").Value = "Exceptional items"
$ws.Range("L2").Value = "Exceptional Items"
